$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update W/X/Y/Z values for rows 2-60 ---
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 0.4538824667597043
$ws.Range("Y2").Value = 80
$ws.Range("Z2").Value = "2025-10-29T23:40:33.767520"
$ws.Range("W3").Value = 1
$ws.Range("X3").Value = 0.5616240759128834
$ws.Range("Y3").Value = 123
$ws.Range("Z3").Value = "2025-10-29T23:40:33.767520"
$ws.Range("W4").Value = 3
$ws.Range("X4").Value = 0.4329311706285884
$ws.Range("Y4").Value = 108
$ws.Range("Z4").Value = "2025-10-29T23:40:33.767520"
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 0.4153959819657586
$ws.Range("Y5").Value = 6
$ws.Range("Z5").Value = "2025-10-29T23:40:33.767520"
$ws.Range("Z6").Value = "2025-10-29T23:40:33.767520"
$ws.Range("Z7").Value = "2025-10-29T23:40:33.767520"
$ws.Range("W8").Value = 6
$ws.Range("X8").Value = 0.5079682182603347
$ws.Range("Y8").Value = 11
$ws.Range("Z8").Value = "2025-10-29T23:40:33.768518"
$ws.Range("Z9").Value = "2025-10-29T23:40:33.768518"
$ws.Range("Z10").Value = "2025-10-29T23:40:33.768518"
$ws.Range("Z11").Value = "2025-10-29T23:40:33.768518"
$ws.Range("W12").Value = 4
$ws.Range("X12").Value = 0.5614880310328125
$ws.Range("Y12").Value = 89
$ws.Range("Z12").Value = "2025-10-29T23:40:33.768518"
$ws.Range("Z13").Value = "2025-10-29T23:40:33.768518"
$ws.Range("W14").Value = 1
$ws.Range("X14").Value = 0.4062858371373469
$ws.Range("Y14").Value = 151
$ws.Range("Z14").Value = "2025-10-29T23:40:33.769517"
$ws.Range("Z15").Value = "2025-10-29T23:40:33.769517"
$ws.Range("W16").Value = 4
$ws.Range("X16").Value = 0.405083825348819
$ws.Range("Y16").Value = 21
$ws.Range("Z16").Value = "2025-10-29T23:40:33.769517"
$ws.Range("Z17").Value = "2025-10-29T23:40:33.769517"
$ws.Range("Z18").Value = "2025-10-29T23:40:33.769517"
$ws.Range("Z19").Value = "2025-10-29T23:40:33.769517"
$ws.Range("Z20").Value = "2025-10-29T23:40:33.770517"
$ws.Range("Z21").Value = "2025-10-29T23:40:33.770517"
$ws.Range("Z22").Value = "2025-10-29T23:40:33.770517"
$ws.Range("Z23").Value = "2025-10-29T23:40:33.770517"
$ws.Range("Z24").Value = "2025-10-29T23:40:33.770517"
$ws.Range("Z25").Value = "2025-10-29T23:40:33.770517"
$ws.Range("W26").Value = 3
$ws.Range("X26").Value = 0.4705137712668338
$ws.Range("Y26").Value = 135
$ws.Range("Z26").Value = "2025-10-29T23:40:33.771517"
$ws.Range("W27").Value = 3
$ws.Range("X27").Value = 0.5541934359909122
$ws.Range("Y27").Value = 145
$ws.Range("Z27").Value = "2025-10-29T23:40:33.771517"
$ws.Range("W28").Value = 2
$ws.Range("X28").Value = 0.4969659942717967
$ws.Range("Y28").Value = 134
$ws.Range("Z28").Value = "2025-10-29T23:40:33.771517"
$ws.Range("W29").Value = 4
$ws.Range("X29").Value = 0.4636006949943728
$ws.Range("Y29").Value = 31
$ws.Range("Z29").Value = "2025-10-29T23:40:33.771517"
$ws.Range("W30").Value = 9
$ws.Range("X30").Value = 0.4641560129943472
$ws.Range("Y30").Value = 108
$ws.Range("Z30").Value = "2025-10-29T23:40:33.797081"
$ws.Range("W31").Value = 7
$ws.Range("X31").Value = 0.4978905520555126
$ws.Range("Y31").Value = 22
$ws.Range("Z31").Value = "2025-10-29T23:40:33.797081"
$ws.Range("W32").Value = 9
$ws.Range("X32").Value = 0.5290345580818899
$ws.Range("Y32").Value = 11
$ws.Range("Z32").Value = "2025-10-29T23:40:33.797081"
$ws.Range("W33").Value = 7
$ws.Range("X33").Value = 0.5381875476204931
$ws.Range("Y33").Value = 6
$ws.Range("Z33").Value = "2025-10-29T23:40:33.797081"
$ws.Range("W34").Value = 5
$ws.Range("X34").Value = 0.5966846281789686
$ws.Range("Y34").Value = 15
$ws.Range("Z34").Value = "2025-10-29T23:40:33.797081"
$ws.Range("Z35").Value = "2025-10-29T23:40:33.797081"
$ws.Range("W36").Value = 9
$ws.Range("X36").Value = 0.4081550283109528
$ws.Range("Y36").Value = 95
$ws.Range("Z36").Value = "2025-10-29T23:40:33.797081"
$ws.Range("Z37").Value = "2025-10-29T23:40:33.797081"
$ws.Range("W38").Value = 13
$ws.Range("X38").Value = 0.4739308912122809
$ws.Range("Y38").Value = 6
$ws.Range("Z38").Value = "2025-10-29T23:40:33.797081"
$ws.Range("Z39").Value = "2025-10-29T23:40:33.798082"
$ws.Range("Z40").Value = "2025-10-29T23:40:33.798082"
$ws.Range("W41").Value = 7
$ws.Range("X41").Value = 0.4727259204758588
$ws.Range("Y41").Value = 79
$ws.Range("Z41").Value = "2025-10-29T23:40:33.799664"
$ws.Range("Z42").Value = "2025-10-29T23:40:33.799664"
$ws.Range("Z43").Value = "2025-10-29T23:40:33.828420"
$ws.Range("W44").Value = 4
$ws.Range("X44").Value = 0.4822074026636463
$ws.Range("Y44").Value = 16
$ws.Range("Z44").Value = "2025-10-29T23:40:33.828420"
$ws.Range("W45").Value = 4
$ws.Range("X45").Value = 0.42961738599068
$ws.Range("Y45").Value = 116
$ws.Range("Z45").Value = "2025-10-29T23:40:33.828420"
$ws.Range("W46").Value = 5
$ws.Range("X46").Value = 0.5721461166512687
$ws.Range("Y46").Value = 16
$ws.Range("Z46").Value = "2025-10-29T23:40:33.828420"
$ws.Range("Z47").Value = "2025-10-29T23:40:33.828420"
$ws.Range("Z48").Value = "2025-10-29T23:40:33.828420"
$ws.Range("W49").Value = 9
$ws.Range("X49").Value = 0.4690142496053366
$ws.Range("Y49").Value = 115
$ws.Range("Z49").Value = "2025-10-29T23:40:33.828420"
$ws.Range("W50").Value = 4
$ws.Range("X50").Value = 0.453356202855057
$ws.Range("Y50").Value = 52
$ws.Range("Z50").Value = "2025-10-29T23:40:33.828420"
$ws.Range("W51").Value = 6
$ws.Range("X51").Value = 0.4444215620941461
$ws.Range("Y51").Value = 17
$ws.Range("Z51").Value = "2025-10-29T23:40:33.828420"
$ws.Range("W52").Value = 5
$ws.Range("X52").Value = 0.5792182599846987
$ws.Range("Y52").Value = 121
$ws.Range("Z52").Value = "2025-10-29T23:40:33.829424"
$ws.Range("Z53").Value = "2025-10-29T23:40:33.829424"
$ws.Range("Z54").Value = "2025-10-29T23:40:33.829424"
$ws.Range("Z55").Value = "2025-10-29T23:40:33.829424"
$ws.Range("Z56").Value = "2025-10-29T23:40:33.829424"
$ws.Range("W57").Value = 12
$ws.Range("X57").Value = 0.5392608545679577
$ws.Range("Y57").Value = 118
$ws.Range("Z57").Value = "2025-10-29T23:40:33.829424"
$ws.Range("Z58").Value = "2025-10-29T23:40:33.829424"
$ws.Range("W59").Value = 9
$ws.Range("X59").Value = 0.450783082786869
$ws.Range("Y59").Value = 144
$ws.Range("Z59").Value = "2025-10-29T23:40:33.829424"
$ws.Range("W60").Value = 6
$ws.Range("X60").Value = 0.4646405864041511
$ws.Range("Y60").Value = 32
$ws.Range("Z60").Value = "2025-10-29T23:40:33.829424"

# --- Update K column for rows 34, 38, 42 ---
$ws.Range("K34").Value = 0.5714285714285714
$ws.Range("K38").Value = 0.4285714285714285
$ws.Range("K42").Value = 0.5714285714285714

# --- Delete rows 61-68 (removed samples) ---
$ws.Range("A61:A68").EntireRow.Delete()

# --- Update conditional formatting range to match new dimension ---
$cf = $ws.Range("A2:Z68").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("A2:Z60"))
